# Commit: "Remember to Fix indexes, include dealer PO"
#
# The sample data in row 5 (second PO block) reused the same "ZEXtest"
# label as row 2, which is incorrect - each PO entry needs its own
# distinct reference. Replace the value in A5 with a new, unique label
# ("ZEXtest3"), which introduces a new shared string, and leave the
# cursor on the cell the user navigated to afterwards (G11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ZEXtest3"

[void]$ws.Range("G11").Select()
